$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H28").Value = 2010.6875
$ws.Range("I28").Value = 2064.5454
$ws.Range("K28").Value = 2064.5454
$ws.Range("M28").Value = -1579.5454
$ws.Range("H51").Value = 2615.5173
$ws.Range("I51").Value = 2450
$ws.Range("J51").Value = 3250
$ws.Range("K51").Value = 2450
$ws.Range("L51").Value = 3250
$ws.Range("M51").Value = -1966
$ws.Range("N51").Value = -4218
$ws.Range("H74").Value = 69023.17999999999
$ws.Range("I74").Value = 148713.58
$ws.Range("K74").Value = 148713.58
$ws.Range("M74").Value = -147777.58
$ws.Range("H77").Value = 69023.17999999999
$ws.Range("I77").Value = 148713.58
$ws.Range("K77").Value = 743567.8999999999
$ws.Range("M77").Value = -738887.8999999999
$ws.Range("H80").Value = 1241.5333
$ws.Range("J80").Value = 1259.4445
$ws.Range("L80").Value = 3778.3335
$ws.Range("N80").Value = -5774.333500000001
$ws.Range("H83").Value = 1241.5333
$ws.Range("J83").Value = 1259.4445
$ws.Range("L83").Value = 11335.0005
$ws.Range("N83").Value = -21319.0005
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H141").Value = 111208.6
$ws.Range("I141").Value = 56898.445
$ws.Range("K141").Value = 170695.335
$ws.Range("M141").Value = -165515.335

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 16481.957
$ws.Range("I45").Value = 14304.105
$ws.Range("J45").Value = 26826.75
$ws.Range("K45").Value = 14304.105
$ws.Range("L45").Value = 26826.75
$ws.Range("M45").Value = -13927.105
$ws.Range("N45").Value = -27580.75
$ws.Range("H61").Value = 3712.25
$ws.Range("I61").Value = 3170.3333
$ws.Range("K61").Value = 3170.3333
$ws.Range("M61").Value = -2958.3333
$ws.Range("H132").Value = 3683.0217
$ws.Range("I132").Value = 3417.1892
$ws.Range("K132").Value = 10251.5676
$ws.Range("M132").Value = -7721.567599999998
$ws.Range("H136").Value = 3712.25
$ws.Range("I136").Value = 3170.3333
$ws.Range("K136").Value = 9510.999899999999
$ws.Range("M136").Value = -6960.999899999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3516.4062
$ws.Range("I86").Value = 3502.6155
$ws.Range("K86").Value = 3502.6155
$ws.Range("M86").Value = -2379.6155
$ws.Range("H89").Value = 3516.4062
$ws.Range("I89").Value = 3502.6155
$ws.Range("K89").Value = 17513.0775
$ws.Range("M89").Value = -11897.0775
$ws.Range("H107").Value = 1945.5
$ws.Range("I107").Value = 1553.05
$ws.Range("J107").Value = 2599.5833
$ws.Range("K107").Value = 1553.05
$ws.Range("L107").Value = 2599.5833
$ws.Range("M107").Value = 366.95
$ws.Range("N107").Value = -6439.5833
$ws.Range("H134").Value = 2069.8928
$ws.Range("I134").Value = 1842.28
$ws.Range("K134").Value = 5526.84
$ws.Range("M134").Value = -2991.84

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27483728
$ws.Range("I31").Value = 4692289
$ws.Range("J31").Value = 71438650
$ws.Range("K31").Value = 4692289
$ws.Range("L31").Value = 71438650
$ws.Range("M31").Value = -4691994
$ws.Range("N31").Value = -71439240
$ws.Range("H34").Value = 27483728
$ws.Range("I34").Value = 4692289
$ws.Range("J34").Value = 71438650
$ws.Range("K34").Value = 4692289
$ws.Range("L34").Value = 71438650
$ws.Range("M34").Value = -4692087
$ws.Range("N34").Value = -71439054
$ws.Range("H105").Value = 3047.25
$ws.Range("J105").Value = 1500
$ws.Range("L105").Value = 1500
$ws.Range("N105").Value = -4994
$ws.Range("H107").Value = 769.5454999999999
$ws.Range("I107").Value = 844.4737
$ws.Range("K107").Value = 844.4737
$ws.Range("M107").Value = 1075.5263
$ws.Range("H120").Value = 250000
$ws.Range("J120").Value = 250000
$ws.Range("L120").Value = 250000
$ws.Range("N120").Value = -257258
$ws.Range("H141").Value = 69125
$ws.Range("J141").Value = 69125
$ws.Range("L141").Value = 69125
$ws.Range("N141").Value = -79485

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 243
$ws.Range("I17").Value = 269.57144
$ws.Range("K17").Value = 808.71432
$ws.Range("M17").Value = -639.71432
$ws.Range("H32").Value = 875
$ws.Range("J32").Value = 900
$ws.Range("L32").Value = 2700
$ws.Range("N32").Value = -3266
$ws.Range("H33").Value = 1912.3334
$ws.Range("I33").Value = 100
$ws.Range("J33").Value = 2818.5
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 16911
$ws.Range("M33").Value = -317
$ws.Range("N33").Value = -17477
$ws.Range("H93").Value = 5954.5
$ws.Range("I93").Value = 3999.5
$ws.Range("J93").Value = 6606.1665
$ws.Range("K93").Value = 11998.5
$ws.Range("L93").Value = 19818.4995
$ws.Range("M93").Value = -10126.5
$ws.Range("N93").Value = -23562.4995
$ws.Range("H98").Value = 17.5
$ws.Range("I98").Value = 25
$ws.Range("K98").Value = 75
$ws.Range("M98").Value = 1423
$ws.Range("H99").Value = 6940.6665
$ws.Range("I99").Value = 1881.6666
$ws.Range("K99").Value = 5644.9998
$ws.Range("M99").Value = -3398.9998
$ws.Range("H107").Value = 17545710
$ws.Range("I107").Value = 111111700
$ws.Range("J107").Value = 2087.5
$ws.Range("K107").Value = 333335100
$ws.Range("L107").Value = 6262.5
$ws.Range("M107").Value = -333333180
$ws.Range("N107").Value = -10102.5
$ws.Range("H134").Value = 6919.8887
$ws.Range("J134").Value = 16000
$ws.Range("L134").Value = 48000
$ws.Range("N134").Value = -58140

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 111111110
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H36").Value = 4820.143
$ws.Range("I36").Value = 5112.5
$ws.Range("J36").Value = 4430.3335
$ws.Range("K36").Value = 5112.5
$ws.Range("L36").Value = 4430.3335
$ws.Range("M36").Value = -4627.5
$ws.Range("N36").Value = -5400.3335
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9537
$ws.Range("H80").Value = 2689.9
$ws.Range("I80").Value = 2655.4443
$ws.Range("K80").Value = 2655.4443
$ws.Range("M80").Value = -1657.4443
$ws.Range("H83").Value = 2689.9
$ws.Range("I83").Value = 2655.4443
$ws.Range("K83").Value = 13277.2215
$ws.Range("M83").Value = -8285.2215
$ws.Range("H97").Value = 2451.5715
$ws.Range("I97").Value = 2106.6843
$ws.Range("K97").Value = 2106.6843
$ws.Range("M97").Value = -1610.6843
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 29289.916
$ws.Range("I126").Value = 34147.9
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 102443.7
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -99973.70000000001
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 13302.066
$ws.Range("I132").Value = 14938.41
$ws.Range("J132").Value = 2665.8333
$ws.Range("K132").Value = 44815.23
$ws.Range("L132").Value = 7997.499899999999
$ws.Range("M132").Value = -42285.23
$ws.Range("N132").Value = -13057.4999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1610.6957
$ws.Range("I22").Value = 1575.75
$ws.Range("J22").Value = 1690.5714
$ws.Range("K22").Value = 1575.75
$ws.Range("L22").Value = 1690.5714
$ws.Range("M22").Value = -1280.75
$ws.Range("N22").Value = -2280.5714
$ws.Range("H23").Value = 10508333
$ws.Range("I23").Value = 683333
$ws.Range("J23").Value = 20333334
$ws.Range("K23").Value = 683333
$ws.Range("L23").Value = 20333334
$ws.Range("M23").Value = -683103
$ws.Range("N23").Value = -20333794
$ws.Range("H27").Value = 1610.6957
$ws.Range("I27").Value = 1575.75
$ws.Range("J27").Value = 1690.5714
$ws.Range("K27").Value = 1575.75
$ws.Range("L27").Value = 1690.5714
$ws.Range("M27").Value = -1468.75
$ws.Range("N27").Value = -1904.5714
$ws.Range("H38").Value = 100000
$ws.Range("I38").Value = 100000
$ws.Range("K38").Value = 100000
$ws.Range("M38").Value = -99590
$ws.Range("H68").Value = 7999
$ws.Range("I68").Value = 7999
$ws.Range("K68").Value = 7999
$ws.Range("M68").Value = -7250
$ws.Range("H71").Value = 7999
$ws.Range("I71").Value = 7999
$ws.Range("K71").Value = 39995
$ws.Range("M71").Value = -36251
$ws.Range("H122").Value = 594916.0600000001
$ws.Range("I122").Value = 1004876.6
$ws.Range("K122").Value = 3014629.8
$ws.Range("M122").Value = -3012179.8
$ws.Range("H132").Value = 4213.7617
$ws.Range("I132").Value = 3422.9412
$ws.Range("J132").Value = 7574.75
$ws.Range("K132").Value = 10268.8236
$ws.Range("L132").Value = 22724.25
$ws.Range("M132").Value = -7738.8236
$ws.Range("N132").Value = -27784.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1259.7142
$ws.Range("I107").Value = 1012
$ws.Range("K107").Value = 3036
$ws.Range("M107").Value = -1116
$ws.Range("H132").Value = 3031.7097
$ws.Range("I132").Value = 3099.5186
$ws.Range("M132").Value = -6768.5558
$ws.Range("H136").Value = 2968
$ws.Range("J136").Value = 3875
$ws.Range("L136").Value = 11625
$ws.Range("N136").Value = -16725
